$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hospitals")

# A new hospital ("Royal Hospital for Children and Young People") was added to the
# Scotland hospital list. Insert a fresh row at position 34 (this pushes the existing
# rows 34-42 down to 35-43, exactly like a manual row insert in Excel).
$ws.Rows.Item(34).Insert()

# The Golden Jubilee National Hospital's NHS Board label was also corrected from
# "National" to "Golden Jubilee Foundation" (now on row 37 after the insert above).
$ws.Range("A37").Value = "Golden Jubilee Foundation"

# Fill in the details for the newly added hospital row.
$ws.Range("A34").Value = "Lothian"
$ws.Range("C34").Value = "Royal Hospital for Children and Young People"
$ws.Range("D34").Value = "50 Little France Crescent, Edinburgh, EH16 4TJ"
$ws.Range("B34").Value = "S319H"
$ws.Range("E34").Value = 73
$ws.Range("F34").Value = "EH16 4TJ"

# Re-fit the columns now that some of the new text is wider than the previous content,
# and leave the selection where the editor's cursor ended up.
for ($i = 1; $i -le 22; $i++) {
    $ws.Columns.Item($i).AutoFit()
}

$ws.Activate()
$ws.Range("E34").Select()
